$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")

# Update counts (F column) for events whose want-to-go count changed
$ws.Cells.Item(2,6).Value = 8
$ws.Cells.Item(3,6).Value = 169
$ws.Cells.Item(4,6).Value = 80
$ws.Cells.Item(6,6).Value = 539
$ws.Cells.Item(7,6).Value = 1660
$ws.Cells.Item(8,6).Value = 16
$ws.Cells.Item(11,6).Value = 1572
$ws.Cells.Item(12,6).Value = 125
$ws.Cells.Item(13,6).Value = 54
$ws.Cells.Item(14,6).Value = 389
$ws.Cells.Item(15,6).Value = 258
$ws.Cells.Item(17,6).Value = 8
$ws.Cells.Item(19,6).Value = 27

# Insert new row for 江西·次元星河国风动漫游戏嘉年华 (2024-07-27)
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).ClearFormats()
$ws.Cells.Item(20,1).Copy()
$ws.Cells.Item(21,1).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(21,2).Value = "'2024-07-27"
$ws.Cells.Item(21,3).Value = "江西·次元星河国风动漫游戏嘉年华"
$ws.Cells.Item(21,4).Value = "九龙大道1177号 南昌绿地国际博览中心"
$ws.Cells.Item(21,5).Value = "2024.07.27 10:00-07.28 17:00"
$ws.Cells.Item(21,6).Value = 37
$ws.Cells.Item(21,7).Value = 39.9
$ws.Cells.Item(21,8).Value = "https://show.bilibili.com/platform/detail.html?id=85493"
$ws.Cells.Item(21,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/HJ7TF5zx1714367786872.jpeg"

# Re-number the index column (A = row-1) for the new row and every shifted row
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(25,1).Value = 24

# Fix up F values on rows shifted down by the insert
$ws.Cells.Item(22,6).Value = 280
$ws.Cells.Item(24,6).Value = 220
$ws.Cells.Item(25,6).Value = 216

$ws = $wb.Worksheets.Item("全部类型")

# Update counts (F column) for events whose want-to-go count changed
$ws.Cells.Item(2,6).Value = 8
$ws.Cells.Item(3,6).Value = 169
$ws.Cells.Item(4,6).Value = 80
$ws.Cells.Item(6,6).Value = 539
$ws.Cells.Item(7,6).Value = 1660
$ws.Cells.Item(9,6).Value = 16
$ws.Cells.Item(12,6).Value = 1572
$ws.Cells.Item(13,6).Value = 125
$ws.Cells.Item(14,6).Value = 54
$ws.Cells.Item(15,6).Value = 389
$ws.Cells.Item(16,6).Value = 258
$ws.Cells.Item(18,6).Value = 8
$ws.Cells.Item(20,6).Value = 27

# Insert new row for 江西·次元星河国风动漫游戏嘉年华 (2024-07-27)
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).ClearFormats()
$ws.Cells.Item(21,1).Copy()
$ws.Cells.Item(22,1).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(22,2).Value = "'2024-07-27"
$ws.Cells.Item(22,3).Value = "江西·次元星河国风动漫游戏嘉年华"
$ws.Cells.Item(22,4).Value = "九龙大道1177号 南昌绿地国际博览中心"
$ws.Cells.Item(22,5).Value = "2024.07.27 10:00-07.28 17:00"
$ws.Cells.Item(22,6).Value = 37
$ws.Cells.Item(22,7).Value = 39.9
$ws.Cells.Item(22,8).Value = "https://show.bilibili.com/platform/detail.html?id=85493"
$ws.Cells.Item(22,9).Value = "//i2.hdslb.com/bfs/openplatform/202404/HJ7TF5zx1714367786872.jpeg"

# Re-number the index column (A = row-1) for the new row and every shifted row
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(25,1).Value = 24
$ws.Cells.Item(26,1).Value = 25

# Fix up F values on rows shifted down by the insert
$ws.Cells.Item(23,6).Value = 280
$ws.Cells.Item(25,6).Value = 220
$ws.Cells.Item(26,6).Value = 216

